# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2  = 3
    3  = 6
    4  = 3
    5  = 5
    6  = 4
    7  = 10
    8  = 4
    9  = 5
    10 = 2
    11 = 2
    12 = 8
    13 = 5
    14 = 3
    15 = 10
    16 = 6
    17 = 5
    18 = 4
    19 = 6
    20 = 10
    21 = 3
    22 = 3
    23 = 10
    24 = 6
    25 = 6
    26 = 8
    27 = 5
    28 = 2
    29 = 6
    30 = 4
    31 = 0
    32 = 6
    33 = 7
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
